# edit.ps1 - applies the "count2: adding CD4 count result in lab test" change
# to the lab.xlsx XLSForm workbook (survey / choices / settings sheets).

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# -----------------------------------------------------------------
# 1) SURVEY sheet
#    Insert a new "note" row right after the existing "CD4 count
#    result" row (row 25) that tells the user to schedule a new CD4
#    lab count, shown only when the count result needs follow up.
#    The row that used to follow (select_one load / Viral Load) gets
#    pushed down and its "relevant" column is tightened so it only
#    shows up for viral-load tests.
# -----------------------------------------------------------------

$survey.Rows.Item(26).Insert()

$survey.Range("A26").Value = "select_one snooze"
$survey.Range("B26").Value = "this"
$survey.Range("C26").Value = "This patient should be scheduled for a CD4 Lab count based on their result and date of last lab appointment. "
$survey.Range("D26").Value = "${count} = 'unstable' or ${count} = 'stable' or ${count} = 'inconclusive' or ${count} = 'unknown'"

# Row 27 is now the old "Viral Load" row (select_one load / load /
# "Viral Load:") - its relevant condition grows an extra clause.
$survey.Range("D27").Value = "${result} = 'yes' and ${test} = 'viral'"

# Column width tweaks (label/relevant columns got wider to fit the
# new, longer text).
$survey.Columns.Item(3).ColumnWidth = 80.66666666666667
$survey.Columns.Item(4).ColumnWidth = 71.66666666666667

# -----------------------------------------------------------------
# 2) CHOICES sheet
#    The old two-option "count" list (unsuppressed/suppressed) is
#    replaced by a richer four-option list (stable/unstable/
#    inconclusive/unknown). Two brand-new "snooze" choices are added
#    at the bottom for the new note/snooze question.
# -----------------------------------------------------------------

# Insert two extra rows so the 2-row block becomes a 4-row block.
$choices.Range("7:8").Insert()

$choices.Range("A7").Value = "count"
$choices.Range("B7").Value = "stable"
$choices.Range("C7").Value = "Stable or Suppressed (CD4 count is 350 or above 350 cells/mm3)"

$choices.Range("A8").Value = "count"
$choices.Range("B8").Value = "unstable"
$choices.Range("C8").Value = "Unstable or Unsuppressed `n(CD4 count is below 350 cells/mm3)`n"
# The embedded line breaks make the engine auto-grow this row; put the
# height back the way it was (matches every other data row).
$choices.Rows.Item(8).AutoFit()

$choices.Range("A9").Value = "count"
$choices.Range("B9").Value = "inconclusive"
$choices.Range("C9").Value = "Inconclusive (Lab test did not work, so patient will need to have a new lab test appointment)"

$choices.Range("A10").Value = "count"
$choices.Range("B10").Value = "unknown"
$choices.Range("C10").Value = "Unknown / no result (Lab test did not give a result. Close this Task and set up new lab test appointment for patient)"

# Two new "snooze" choices appended at the bottom of the sheet.
$choices.Range("A13").Value = "snooze"
$choices.Range("B13").Value = "snooze1"
$choices.Range("C13").Value = "Okay,I will schedule an appointment for a lab visit"

$choices.Range("A14").Value = "snooze"
$choices.Range("B14").Value = "snooze2"
$choices.Range("C14").Value = "Remind me to schedule an appointment in 2 days "

# Column width tweak (label column got wider to fit the new text).
$choices.Columns.Item(3).ColumnWidth = 83.83333333333333
